$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 9.454274208217237
$ws.Cells.Item(2, 4).Value = 8.067803511874837
$ws.Cells.Item(2, 5).Value = 13.1591483127652
$ws.Cells.Item(2, 6).Value = 37.88371238076034
$ws.Cells.Item(2, 7).Value = 3.695655277174606
$ws.Cells.Item(2, 9).Value = 29.47010556688408
$ws.Cells.Item(2, 10).Value = 10.33558898771845
$ws.Cells.Item(2, 11).Value = 16.84455902426148
$ws.Cells.Item(2, 12).Value = 10.30065057450186
$ws.Cells.Item(2, 13).Value = 18.81103664535916
$ws.Cells.Item(2, 15).Value = 29.23913613188876

$ws.Cells.Item(3, 3).Value = 9.422657315158791
$ws.Cells.Item(3, 4).Value = 8.041219994027994
$ws.Cells.Item(3, 5).Value = 13.18111054357068
$ws.Cells.Item(3, 6).Value = 38.05170069082426
$ws.Cells.Item(3, 7).Value = 3.69780740067419
$ws.Cells.Item(3, 9).Value = 29.60235104795779
$ws.Cells.Item(3, 10).Value = 10.36104859610437
$ws.Cells.Item(3, 11).Value = 16.37389424061714
$ws.Cells.Item(3, 12).Value = 10.32311689373139
$ws.Cells.Item(3, 13).Value = 18.62133038279577
$ws.Cells.Item(3, 15).Value = 29.38924300627838

$ws.Cells.Item(4, 3).Value = 9.404468092232813
$ws.Cells.Item(4, 4).Value = 8.025599766373826
$ws.Cells.Item(4, 5).Value = 13.19612939360933
$ws.Cells.Item(4, 6).Value = 38.16391510312033
$ws.Cells.Item(4, 7).Value = 3.699198436834952
$ws.Cells.Item(4, 9).Value = 29.69001307400459
$ws.Cells.Item(4, 10).Value = 10.3775509716988
$ws.Cells.Item(4, 11).Value = 16.07847783679615
$ws.Cells.Item(4, 12).Value = 10.3377329332411
$ws.Cells.Item(4, 13).Value = 18.50534300632922
$ws.Cells.Item(4, 15).Value = 29.48799108588045

$ws.Cells.Item(5, 3).Value = 9.397368446766412
$ws.Cells.Item(5, 4).Value = 8.019414562140458
$ws.Cells.Item(5, 5).Value = 13.20263586212579
$ws.Cells.Item(5, 6).Value = 38.21191876232358
$ws.Cells.Item(5, 7).Value = 3.69978285817632
$ws.Cells.Item(5, 9).Value = 29.72735883232113
$ws.Cells.Item(5, 10).Value = 10.38449518328983
$ws.Cells.Item(5, 11).Value = 15.95663541290967
$ws.Cells.Item(5, 12).Value = 10.34389620369597
$ws.Cells.Item(5, 13).Value = 18.45824036454642
$ws.Cells.Item(5, 15).Value = 29.52988488610177

$ws.Cells.Item(6, 3).Value = 9.396208574761067
$ws.Cells.Item(6, 4).Value = 8.018398497597602
$ws.Cells.Item(6, 5).Value = 13.2037395922265
$ws.Cells.Item(6, 6).Value = 38.22002701133964
$ws.Cells.Item(6, 7).Value = 3.699880963274019
$ws.Cells.Item(6, 9).Value = 29.73365801022145
$ws.Cells.Item(6, 10).Value = 10.38566152997483
$ws.Cells.Item(6, 11).Value = 15.93632035418446
$ws.Cells.Item(6, 12).Value = 10.34493213437554
$ws.Cells.Item(6, 13).Value = 18.45042999974259
$ws.Cells.Item(6, 15).Value = 29.53694111185797

$ws.Cells.Item(7, 3).Value = 9.404371071902997
$ws.Cells.Item(7, 4).Value = 8.025515616406407
$ws.Cells.Item(7, 5).Value = 13.19621557797813
$ws.Cells.Item(7, 6).Value = 38.16455329111349
$ws.Cells.Item(7, 7).Value = 3.699206247361001
$ws.Cells.Item(7, 9).Value = 29.6905101645315
$ws.Cells.Item(7, 10).Value = 10.37764373482631
$ws.Cells.Item(7, 11).Value = 16.07684031312841
$ws.Cells.Item(7, 12).Value = 10.33781521397947
$ws.Cells.Item(7, 13).Value = 18.50470705150168
$ws.Cells.Item(7, 15).Value = 29.4885493892521

$ws.Cells.Item(8, 3).Value = 9.443121755010536
$ws.Cells.Item(8, 4).Value = 8.05849398288356
$ws.Cells.Item(8, 5).Value = 13.16640278800813
$ws.Cells.Item(8, 6).Value = 37.93974946217467
$ws.Cells.Item(8, 7).Value = 3.696382912378625
$ws.Cells.Item(8, 9).Value = 29.51436105366991
$ws.Cells.Item(8, 10).Value = 10.34418722222809
$ws.Cells.Item(8, 11).Value = 16.68369906026793
$ws.Cells.Item(8, 12).Value = 10.30822673978236
$ws.Cells.Item(8, 13).Value = 18.74554717403504
$ws.Cells.Item(8, 15).Value = 29.28952576433245

$ws.Cells.Item(9, 3).Value = 9.52859420960603
$ws.Cells.Item(9, 4).Value = 8.12857255527507
$ws.Cells.Item(9, 5).Value = 13.12009155820867
$ws.Cells.Item(9, 6).Value = 37.5710925070282
$ws.Cells.Item(9, 7).Value = 3.691396261688709
$ws.Cells.Item(9, 9).Value = 29.2203169326485
$ws.Cells.Item(9, 10).Value = 10.28545660038914
$ws.Cells.Item(9, 11).Value = 17.81637344599734
$ws.Cells.Item(9, 12).Value = 10.25669992510703
$ws.Cells.Item(9, 13).Value = 19.22001890188118
$ws.Cells.Item(9, 15).Value = 28.9515509856327

$ws.Cells.Item(10, 3).Value = 9.596851976491399
$ws.Cells.Item(10, 4).Value = 8.183130545843053
$ws.Cells.Item(10, 5).Value = 13.09344770896012
$ws.Cells.Item(10, 6).Value = 37.34454938378286
$ws.Cells.Item(10, 7).Value = 3.688064253536246
$ws.Cells.Item(10, 9).Value = 29.03575052801522
$ws.Cells.Item(10, 10).Value = 10.24646331118232
$ws.Cells.Item(10, 11).Value = 18.60582251027675
$ws.Cells.Item(10, 12).Value = 10.22277113804436
$ws.Cells.Item(10, 13).Value = 19.56763627584622
$ws.Cells.Item(10, 15).Value = 28.73524942512333

$ws.Cells.Item(11, 3).Value = 9.629016828984598
$ws.Cells.Item(11, 4).Value = 8.208566655188486
$ws.Cells.Item(11, 5).Value = 13.08292333748741
$ws.Cells.Item(11, 6).Value = 37.25117481068214
$ws.Cells.Item(11, 7).Value = 3.686619702976558
$ws.Cells.Item(11, 9).Value = 28.95865051907206
$ws.Cells.Item(11, 10).Value = 10.22961894794655
$ws.Cells.Item(11, 11).Value = 18.95427205052321
$ws.Cells.Item(11, 12).Value = 10.2081821544448
$ws.Cells.Item(11, 13).Value = 19.72507865388893
$ws.Cells.Item(11, 15).Value = 28.64382569919702

$ws.Cells.Item(12, 3).Value = 9.641350450524071
$ws.Cells.Item(12, 4).Value = 8.218282933249274
$ws.Cells.Item(12, 5).Value = 13.07916699449442
$ws.Cells.Item(12, 6).Value = 37.21721374394859
$ws.Cells.Item(12, 7).Value = 3.686082870034479
$ws.Cells.Item(12, 9).Value = 28.93044395071665
$ws.Cells.Item(12, 10).Value = 10.22336838274056
$ws.Cells.Item(12, 11).Value = 19.08458253729797
$ws.Cells.Item(12, 12).Value = 10.20277873774287
$ws.Cells.Item(12, 13).Value = 19.78455855707361
$ws.Cells.Item(12, 15).Value = 28.6102111587948

$ws.Cells.Item(13, 3).Value = 9.638687466933343
$ws.Cells.Item(13, 4).Value = 8.216186689633211
$ws.Cells.Item(13, 5).Value = 13.0799658132416
$ws.Cells.Item(13, 6).Value = 37.22446560671094
$ws.Cells.Item(13, 7).Value = 3.6861980343514
$ws.Cells.Item(13, 9).Value = 28.9364746815981
$ws.Cells.Item(13, 10).Value = 10.22470886769629
$ws.Cells.Item(13, 11).Value = 19.05659233059705
$ws.Cells.Item(13, 12).Value = 10.20393708073787
$ws.Cells.Item(13, 13).Value = 19.77175536888937
$ws.Cells.Item(13, 15).Value = 28.6174058777288

$ws.Cells.Item(14, 3).Value = 9.630028486437435
$ws.Cells.Item(14, 4).Value = 8.2093643589794
$ws.Cells.Item(14, 5).Value = 13.08260971419259
$ws.Cells.Item(14, 6).Value = 37.24835276555687
$ws.Cells.Item(14, 7).Value = 3.686575333535543
$ws.Cells.Item(14, 9).Value = 28.95631009854818
$ws.Cells.Item(14, 10).Value = 10.22910214717977
$ws.Cells.Item(14, 11).Value = 18.96502619977446
$ws.Cells.Item(14, 12).Value = 10.20773518734996
$ws.Cells.Item(14, 13).Value = 19.72997514041783
$ws.Cells.Item(14, 15).Value = 28.64104003907791

$ws.Cells.Item(15, 3).Value = 9.624744397171362
$ws.Cells.Item(15, 4).Value = 8.205196308122138
$ws.Cells.Item(15, 5).Value = 13.08425898814001
$ws.Cells.Item(15, 6).Value = 37.26316655131895
$ws.Cells.Item(15, 7).Value = 3.686807765321206
$ws.Cells.Item(15, 9).Value = 28.96858882097379
$ws.Cells.Item(15, 10).Value = 10.23180981592488
$ws.Cells.Item(15, 11).Value = 18.9087228195656
$ws.Cells.Item(15, 12).Value = 10.21007739671464
$ws.Cells.Item(15, 13).Value = 19.70436409922041
$ws.Cells.Item(15, 15).Value = 28.65564769799041

$ws.Cells.Item(16, 3).Value = 9.5947717765194
$ws.Cells.Item(16, 4).Value = 8.181480268654155
$ws.Cells.Item(16, 5).Value = 13.09416757566902
$ws.Cells.Item(16, 6).Value = 37.35084693102036
$ws.Cells.Item(16, 7).Value = 3.688160086604072
$ws.Cells.Item(16, 9).Value = 29.0409274994415
$ws.Cells.Item(16, 10).Value = 10.24758207342579
$ws.Cells.Item(16, 11).Value = 18.5828265142724
$ws.Cells.Item(16, 12).Value = 10.22374153517012
$ws.Cells.Item(16, 13).Value = 19.55732964599438
$ws.Cells.Item(16, 15).Value = 28.74136474354432

$ws.Cells.Item(17, 3).Value = 9.576665172996938
$ws.Cells.Item(17, 4).Value = 8.167086134972681
$ws.Cells.Item(17, 5).Value = 13.10065461078294
$ws.Cells.Item(17, 6).Value = 37.40711965454124
$ws.Cells.Item(17, 7).Value = 3.689007890681036
$ws.Cells.Item(17, 9).Value = 29.08706415605685
$ws.Cells.Item(17, 10).Value = 10.25748642254398
$ws.Cells.Item(17, 11).Value = 18.38008773310591
$ws.Cells.Item(17, 12).Value = 10.23234024387354
$ws.Cells.Item(17, 13).Value = 19.46692265296042
$ws.Cells.Item(17, 15).Value = 28.79573737891632

$ws.Cells.Item(18, 3).Value = 9.566356008714484
$ws.Cells.Item(18, 4).Value = 8.158865405778958
$ws.Cells.Item(18, 5).Value = 13.10453604115111
$ws.Cells.Item(18, 6).Value = 37.44039712252535
$ws.Cells.Item(18, 7).Value = 3.689502229552263
$ws.Cells.Item(18, 9).Value = 29.11424626981393
$ws.Cells.Item(18, 10).Value = 10.26326730847327
$ws.Cells.Item(18, 11).Value = 18.26247856686475
$ws.Cells.Item(18, 12).Value = 10.23736559266427
$ws.Cells.Item(18, 13).Value = 19.4148603487919
$ws.Cells.Item(18, 15).Value = 28.82766722583505

$ws.Cells.Item(19, 3).Value = 9.562883784932946
$ws.Cells.Item(19, 4).Value = 8.15609217628953
$ws.Cells.Item(19, 5).Value = 13.1058760488355
$ws.Cells.Item(19, 6).Value = 37.45182057087939
$ws.Cells.Item(19, 7).Value = 3.689670757245889
$ws.Cells.Item(19, 9).Value = 29.12356044406799
$ws.Cells.Item(19, 10).Value = 10.26523908836572
$ws.Cells.Item(19, 11).Value = 18.22248987778324
$ws.Cells.Item(19, 12).Value = 10.23908077445575
$ws.Cells.Item(19, 13).Value = 19.39722345202072
$ws.Cells.Item(19, 15).Value = 28.83859072866342

$ws.Cells.Item(20, 3).Value = 9.578581804220738
$ws.Cells.Item(20, 4).Value = 8.168612404378687
$ws.Cells.Item(20, 5).Value = 13.09994850709839
$ws.Cells.Item(20, 6).Value = 37.40103501195801
$ws.Cells.Item(20, 7).Value = 3.688916946980964
$ws.Cells.Item(20, 9).Value = 29.08208599828406
$ws.Cells.Item(20, 10).Value = 10.25642338048131
$ws.Cells.Item(20, 11).Value = 18.40177382602572
$ws.Cells.Item(20, 12).Value = 10.23141666174281
$ws.Cells.Item(20, 13).Value = 19.47655340778675
$ws.Cells.Item(20, 15).Value = 28.78988138941521

$ws.Cells.Item(21, 3).Value = 9.632567727712171
$ws.Cells.Item(21, 4).Value = 8.211365995590695
$ws.Cells.Item(21, 5).Value = 13.08182692530965
$ws.Cells.Item(21, 6).Value = 37.24129854518725
$ws.Cells.Item(21, 7).Value = 3.686464235551158
$ws.Cells.Item(21, 9).Value = 28.95045707584168
$ws.Cells.Item(21, 10).Value = 10.22780826384648
$ws.Cells.Item(21, 11).Value = 18.99196667895927
$ws.Cells.Item(21, 12).Value = 10.20661630792594
$ws.Cells.Item(21, 13).Value = 19.7422511303407
$ws.Cells.Item(21, 15).Value = 28.63407079637529

$ws.Cells.Item(22, 3).Value = 9.668742054904733
$ws.Cells.Item(22, 4).Value = 8.239796371320201
$ws.Cells.Item(22, 5).Value = 13.07131800702586
$ws.Cells.Item(22, 6).Value = 37.14505191763124
$ws.Cells.Item(22, 7).Value = 3.684920601111618
$ws.Cells.Item(22, 9).Value = 28.87019913279062
$ws.Cells.Item(22, 10).Value = 10.20985265279882
$ws.Cells.Item(22, 11).Value = 19.36809247900683
$ws.Cells.Item(22, 12).Value = 10.19111358639091
$ws.Cells.Item(22, 13).Value = 19.91506491906557
$ws.Cells.Item(22, 15).Value = 28.53810262832805

$ws.Cells.Item(23, 3).Value = 9.649355840184565
$ws.Cells.Item(23, 4).Value = 8.224579400929938
$ws.Cells.Item(23, 5).Value = 13.07680487260507
$ws.Cells.Item(23, 6).Value = 37.19567291868777
$ws.Cells.Item(23, 7).Value = 3.685739053691998
$ws.Cells.Item(23, 9).Value = 28.9125054197213
$ws.Cells.Item(23, 10).Value = 10.21936780615858
$ws.Cells.Item(23, 11).Value = 19.16825735648757
$ws.Cells.Item(23, 12).Value = 10.19932325304404
$ws.Cells.Item(23, 13).Value = 19.82292063646762
$ws.Cells.Item(23, 15).Value = 28.58878518361823

$ws.Cells.Item(24, 3).Value = 9.577714981511454
$ws.Cells.Item(24, 4).Value = 8.16792220741954
$ws.Cells.Item(24, 5).Value = 13.10026726323948
$ws.Cells.Item(24, 6).Value = 37.4037829962989
$ws.Cells.Item(24, 7).Value = 3.68895804105244
$ws.Cells.Item(24, 9).Value = 29.08433457542958
$ws.Cells.Item(24, 10).Value = 10.25690371155507
$ws.Cells.Item(24, 11).Value = 18.39197281191773
$ws.Cells.Item(24, 12).Value = 10.23183395827324
$ws.Cells.Item(24, 13).Value = 19.4721996086516
$ws.Cells.Item(24, 15).Value = 28.79252679451959

$ws.Cells.Item(25, 3).Value = 9.504488441658113
$ws.Cells.Item(25, 4).Value = 8.109058232292774
$ws.Cells.Item(25, 5).Value = 13.13132175389831
$ws.Cells.Item(25, 6).Value = 37.66306467809102
$ws.Cells.Item(25, 7).Value = 3.692686780561647
$ws.Cells.Item(25, 9).Value = 29.2943479837462
$ws.Cells.Item(25, 10).Value = 10.30061230991965
$ws.Cells.Item(25, 11).Value = 17.5169415381147
$ws.Cells.Item(25, 12).Value = 10.26994719113318
$ws.Cells.Item(25, 13).Value = 19.09166684294394
$ws.Cells.Item(25, 15).Value = 29.03736998718004
